$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a paragraph (1-based index) whose text contains a substring.
# ---------------------------------------------------------------------------
function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$text*") {
            return $i
        }
    }
    return -1
}

# Word's wdColor values are encoded 0x00BBGGRR (reverse byte order of the
# "RRGGBB" hex we see in OOXML <w:color w:val="2C3E50"/>).
$accentColor = 5258796   # 0x2C3E50 -> BGR 0x503E2C

# ---------------------------------------------------------------------------
# Change 1: professional-summary sentence - neutralize the language.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2)

# ---------------------------------------------------------------------------
# Change 2: "Discovered systematic race coding errors ..." bullet under the
# Siege Analytics role - same neutralization, but "50M" must land in its own
# bold / colored run (matching the styling already used for the other
# highlighted numbers in that bullet).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial machine",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial machine", 2)

$rng = $d.Content
$null = $rng.Find.Execute("50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 1
$rng.Font.Color = $accentColor

# ---------------------------------------------------------------------------
# Change 3: reorder PROFESSIONAL EXPERIENCE entries.
#   - "Data Products Manager - Helm/Murmuration" moves to right before
#     "Software Engineer - Mautinoa Technologies".
#   - "Research Director - PCCC" moves to right after the Mautinoa section
#     (i.e. right before "Software Engineer - Salsa Labs").
# Both sections are rebuilt fresh at their new locations (preserving their
# heading style and any bold/colored runs) and then the original copies are
# deleted.
# ---------------------------------------------------------------------------

# Each paragraph: @{ text = "..."; style = "Heading 3" | "Normal"; bold = $true/$false(optional, applies to a trailing highlighted run) ; boldText = "..." }
function Insert-ParaBefore($anchorIndex, $text, $style) {
    $anchor = $d.Paragraphs($anchorIndex)
    $anchor.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($anchorIndex)
    $newPara.Range.Text = $text
    $newPara.Style = $style
    return $anchorIndex  # the newly created paragraph now sits at this index
}

function Insert-ParaWithTrailingBoldBefore($anchorIndex, $prefixText, $boldText, $style) {
    $anchor = $d.Paragraphs($anchorIndex)
    $anchor.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($anchorIndex)
    $newPara.Range.Text = $prefixText
    $newPara.Style = $style

    $newPara.Range.InsertAfter($boldText)

    $searchRange = $newPara.Range
    $null = $searchRange.Find.Execute($boldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $searchRange.Font.Bold = 1
    $searchRange.Font.Color = $accentColor
}

# --- Insert the "Data Products Manager - Helm/Murmuration" block right
#     before the Mautinoa Technologies heading. -----------------------------
$anchorIdx = Find-ParaIndex("Software Engineer - Mautinoa Technologies")
Insert-ParaBefore $anchorIdx "Data Products Manager - Helm/Murmuration (Austin, TX) | 2021 - 2023" "Heading 3" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Mautinoa Technologies")
Insert-ParaBefore $anchorIdx "Democratic Electoral Technology" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Mautinoa Technologies")
Insert-ParaBefore $anchorIdx "• Led design and implementation of enterprise-scale multi-tenant data warehouse for geo-referenced demographic, econometric, and electoral data" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Mautinoa Technologies")
Insert-ParaBefore $anchorIdx "• Managed engineering team of 11 professionals while setting technical direction for data architecture" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Mautinoa Technologies")
Insert-ParaWithTrailingBoldBefore $anchorIdx "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by " "57%" "Normal"

# --- Insert the "Research Director - PCCC" block right before the Salsa
#     Labs heading (i.e. right after the Mautinoa section). ----------------
$anchorIdx = Find-ParaIndex("Software Engineer - Salsa Labs")
Insert-ParaBefore $anchorIdx "Research Director - PCCC (Washington, DC) | August 2011 - August 2012" "Heading 3" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Salsa Labs")
Insert-ParaBefore $anchorIdx "Political Research & Data Analysis (FLEEM System)" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Salsa Labs")
Insert-ParaBefore $anchorIdx "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Salsa Labs")
Insert-ParaBefore $anchorIdx "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren" "Normal" | Out-Null

$anchorIdx = Find-ParaIndex("Software Engineer - Salsa Labs")
Insert-ParaBefore $anchorIdx "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver" "Normal" | Out-Null

# ---------------------------------------------------------------------------
# Remove the original copies of both blocks (now further down the document,
# still in their old relative order: Praxis Project -> PCCC -> Helm).
# Delete from the bottom paragraph upward within each block so indices of
# not-yet-deleted paragraphs stay valid.
# ---------------------------------------------------------------------------

$oldBlockMarkers = @(
    "• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver",
    "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren",
    "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys",
    "Political Research & Data Analysis (FLEEM System)",
    "Research Director - PCCC (Washington, DC) | August 2011 - August 2012",
    "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%",
    "• Managed engineering team of 11 professionals while setting technical direction for data architecture",
    "• Led design and implementation of enterprise-scale multi-tenant data warehouse for geo-referenced demographic, econometric, and electoral data",
    "Democratic Electoral Technology",
    "Data Products Manager - Helm/Murmuration (Austin, TX) | 2021 - 2023"
)

# There are now two occurrences of each marker (the freshly inserted one and
# the original one further down) except these markers are unique substrings,
# so find the LAST occurrence (the original, still-undeleted copy) each time.
function Find-LastParaIndex($text) {
    $found = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$text*") {
            $found = $i
        }
    }
    return $found
}

foreach ($marker in $oldBlockMarkers) {
    $idx = Find-LastParaIndex($marker)
    if ($idx -ge 1) {
        $d.Paragraphs($idx).Range.Delete()
    }
}

Write-Output "done"
